# Append two trailing spaces to the first sentence, then append a
# red-colored parenthetical note (typed/saved in three chunks, hence
# three separate runs sharing the same red-colored rPr).

$d = $word.ActiveDocument

# 1) "This is a Microsoft word document." -> "This is a Microsoft word document.  "
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs.First
$r = $p1.Range
$textLen = $r.Text.Length
# position right after the trailing double-space, before the paragraph mark
$pos = $r.Start + $textLen - 1

$enDash = [char]0x2013

# 2) First red run: "(This is a change – Ve"
$run1 = $d.Range($pos, $pos)
$run1.InsertAfter("(This is a change " + $enDash + " Ve")
$run1.Font.Color = 255

# 3) Second red run: "rsion for main branch"
$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for main branch")
$run2.Font.Color = 255

# 4) Third red run: ")"
$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 255

Write-Output ("Paragraph 1 now reads: [" + $p1.Range.Text + "]")
